# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 61 (pushing the
# existing rows 61-69 down to 62-70), and the new row is populated with
# its own data (same dimension/categorisation columns as its neighbours,
# new date + price figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 61, shifting rows 61:69 down to 62:70.
$ws.Rows("61:61").Insert()

# Populate the newly inserted row 61 with the new weekly record.
$ws.Range("A61").Value = 4
$ws.Range("B61").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C61").Value = "Los Lagos"
$ws.Range("D61").Value = 45209
$ws.Range("E61").Value = 10
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100101
$ws.Range("H61").Value = "Berries"
$ws.Range("I61").Value = 100101001
$ws.Range("J61").Value = "Arándano (blue)"
$ws.Range("K61").Value = "Sin especificar"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 100
$ws.Range("N61").Value = 13000
$ws.Range("O61").Value = 13000
$ws.Range("P61").Value = 13000
$ws.Range("Q61").Value = "$/bandeja 2 kilos"
$ws.Range("R61").Value = "Provincia de Curicó"
$ws.Range("S61").Value = 6500
$ws.Range("T61").Value = 2
